$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1732.3334
$ws.Range("I12").Value = 433.33334
$ws.Range("J12").Value = 2165.3333
$ws.Range("K12").Value = 433.33334
$ws.Range("L12").Value = 2165.3333
$ws.Range("M12").Value = -263.33334
$ws.Range("N12").Value = -2505.3333
$ws.Range("H34").Value = 1144
$ws.Range("I34").Value = 1144
$ws.Range("K34").Value = 1144
$ws.Range("M34").Value = -941
$ws.Range("H36").Value = 1144
$ws.Range("I36").Value = 1144
$ws.Range("K36").Value = 1144
$ws.Range("M36").Value = -429
$ws.Range("H74").Value = 9516.933999999999
$ws.Range("I74").Value = 8640.909
$ws.Range("J74").Value = 11926
$ws.Range("K74").Value = 8640.909
$ws.Range("L74").Value = 11926
$ws.Range("M74").Value = -7704.909
$ws.Range("N74").Value = -13798
$ws.Range("H77").Value = 9516.933999999999
$ws.Range("I77").Value = 8640.909
$ws.Range("J77").Value = 11926
$ws.Range("K77").Value = 43204.545
$ws.Range("L77").Value = 59630
$ws.Range("M77").Value = -38524.545
$ws.Range("N77").Value = -68990
$ws.Range("H82").Value = 402
$ws.Range("I82").Value = 402
$ws.Range("K82").Value = 1206
$ws.Range("M82").Value = -800
$ws.Range("H85").Value = 402
$ws.Range("I85").Value = 402
$ws.Range("K85").Value = 1206
$ws.Range("M85").Value = 198
$ws.Range("H113").Value = 7245.846
$ws.Range("H116").Value = 11166
$ws.Range("I116").Value = 9400.799999999999
$ws.Range("K116").Value = 9400.799999999999
$ws.Range("M116").Value = -5958.799999999999
$ws.Range("H137").Value = 7033.6665
$ws.Range("I137").Value = 4999.6665
$ws.Range("K137").Value = 14998.9995
$ws.Range("M137").Value = -12448.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5053.0977
$ws.Range("I61").Value = 3478.5625
$ws.Range("K61").Value = 3478.5625
$ws.Range("M61").Value = -3266.5625
$ws.Range("H74").Value = 6947776.5
$ws.Range("I74").Value = 9527513
$ws.Range("K74").Value = 9527513
$ws.Range("M74").Value = -9526639
$ws.Range("H77").Value = 6947776.5
$ws.Range("I77").Value = 9527513
$ws.Range("K77").Value = 47637565
$ws.Range("M77").Value = -47633197
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H110").Value = 3078.5715
$ws.Range("I110").Value = 2578.7083
$ws.Range("K110").Value = 2578.7083
$ws.Range("M110").Value = -533.7082999999998
$ws.Range("H136").Value = 5053.0977
$ws.Range("I136").Value = 3478.5625
$ws.Range("K136").Value = 10435.6875
$ws.Range("M136").Value = -7885.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H20").Value = 1328.4
$ws.Range("I20").Value = 1334.9286
$ws.Range("J20").Value = 1313.1666
$ws.Range("K20").Value = 1334.9286
$ws.Range("L20").Value = 1313.1666
$ws.Range("M20").Value = -1087.9286
$ws.Range("N20").Value = -1807.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 156.15384
$ws.Range("I7").Value = 83.09999999999999
$ws.Range("J7").Value = 399.66666
$ws.Range("K7").Value = 83.09999999999999
$ws.Range("L7").Value = 399.66666
$ws.Range("M7").Value = 29.90000000000001
$ws.Range("N7").Value = -625.66666
$ws.Range("H25").Value = 1013
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 53044.61
$ws.Range("I31").Value = 5014.143
$ws.Range("K31").Value = 5014.143
$ws.Range("M31").Value = -4719.143
$ws.Range("H34").Value = 53044.61
$ws.Range("I34").Value = 5014.143
$ws.Range("K34").Value = 5014.143
$ws.Range("M34").Value = -4812.143
$ws.Range("H124").Value = 44000
$ws.Range("J124").Value = 44000
$ws.Range("L124").Value = 44000
$ws.Range("N124").Value = -48910
$ws.Range("H132").Value = 3446.3684
$ws.Range("I132").Value = 2786.606
$ws.Range("K132").Value = 8359.818000000001
$ws.Range("M132").Value = -5829.818000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31695.947
$ws.Range("I2").Value = 114.545456
$ws.Range("K2").Value = 687.272736
$ws.Range("M2").Value = -574.272736
$ws.Range("H92").Value = 1700.5
$ws.Range("I92").Value = 226
$ws.Range("J92").Value = 2192
$ws.Range("K92").Value = 678
$ws.Range("L92").Value = 6576
$ws.Range("M92").Value = 570
$ws.Range("N92").Value = -9072
$ws.Range("H132").Value = 6333.25
$ws.Range("I132").Value = 5599.8
$ws.Range("K132").Value = 50398.2
$ws.Range("M132").Value = -47868.2
$ws.Range("H133").Value = 4690.8887
$ws.Range("J133").Value = 19030
$ws.Range("L133").Value = 57090
$ws.Range("N133").Value = -67210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 802.53845
$ws.Range("I2").Value = 90.42856999999999
$ws.Range("J2").Value = 1633.3334
$ws.Range("K2").Value = 90.42856999999999
$ws.Range("L2").Value = 1633.3334
$ws.Range("M2").Value = 22.57143000000001
$ws.Range("N2").Value = -1859.3334
$ws.Range("H3").Value = 4011.875
$ws.Range("I3").Value = 2119
$ws.Range("K3").Value = 2119
$ws.Range("M3").Value = -2003
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 10000
$ws.Range("K4").Value = 10000
$ws.Range("M4").Value = -9888
$ws.Range("H70").Value = 7051
$ws.Range("I70").Value = 6953.7144
$ws.Range("J70").Value = 7164.5
$ws.Range("K70").Value = 6953.7144
$ws.Range("L70").Value = 7164.5
$ws.Range("M70").Value = -6683.7144
$ws.Range("N70").Value = -7704.5
$ws.Range("H73").Value = 7051
$ws.Range("I73").Value = 6953.7144
$ws.Range("J73").Value = 7164.5
$ws.Range("K73").Value = 6953.7144
$ws.Range("L73").Value = 7164.5
$ws.Range("M73").Value = -6017.7144
$ws.Range("N73").Value = -9036.5
$ws.Range("H107").Value = 2517.1667
$ws.Range("I107").Value = 788.25
$ws.Range("J107").Value = 5975
$ws.Range("K107").Value = 788.25
$ws.Range("L107").Value = 5975
$ws.Range("M107").Value = 1131.75
$ws.Range("N107").Value = -9815

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3262.2964
$ws.Range("I93").Value = 798.3158
$ws.Range("J93").Value = 9114.25
$ws.Range("K93").Value = 798.3158
$ws.Range("L93").Value = 9114.25
$ws.Range("M93").Value = 449.6842
$ws.Range("N93").Value = -11610.25
$ws.Range("H100").Value = 3422.1765
$ws.Range("I100").Value = 2288.5454
$ws.Range("J100").Value = 5500.5
$ws.Range("K100").Value = 2288.5454
$ws.Range("L100").Value = 5500.5
$ws.Range("M100").Value = -1747.5454
$ws.Range("N100").Value = -6582.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 102
$ws.Range("I7").Value = 102
$ws.Range("K7").Value = 102
$ws.Range("M7").Value = 11
$ws.Range("H113").Value = 362.54544
$ws.Range("I113").Value = 289.47058
$ws.Range("J113").Value = 611
$ws.Range("K113").Value = 868.41174
$ws.Range("L113").Value = 1833
$ws.Range("M113").Value = 1301.58826
$ws.Range("N113").Value = -6173
